$wb = $excel.ActiveWorkbook

# --- Update DomesticFTA sheet selection (tab no longer active/selected) ---
$domestic = $wb.Worksheets.Item("DomesticFTA")
$domestic.Activate() | Out-Null
$domestic.Range("A4:I4").Select() | Out-Null

# --- Add the new InternationalFT sheet after DomesticFTA ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$intl = $wb.Worksheets.Add($null, $lastSheet)
$intl.Name = "InternationalFT"

# Header row
$intl.Range("A1").Value = "DataBinding"
$intl.Range("B1").Value = "ReciverBankName"
$intl.Range("C1").Value = "ReceiverName"
$intl.Range("D1").Value = "ReceAccNum"
$intl.Range("E1").Value = "SwiftMsg"
$intl.Range("F1").Value = "Amount"
$intl.Range("G1").Value = "TransferType"
$intl.Range("H1").Value = "DOT"
$intl.Range("I1").Value = "TransferDesc"

# Data row
$intl.Range("A2").Value = "Data001"
$intl.Range("B2").Value = "RBS"
$intl.Range("C2").Value = "Smith"
$intl.Range("D2").Value = 1234556657
$intl.Range("E2").Value = "MT103"
$intl.Range("F2").Value = 6
$intl.Range("G2").Value = "Domestic Transfer"
$intl.Range("H2").Value = 120481
$intl.Range("I2").Value = "TransferDesc"

# Copy the date number-format from DomesticFTA's H4 (reuses existing style, avoids creating a new numFmt)
$domestic.Range("H4").Copy() | Out-Null
$intl.Range("H2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$intl.Range("H2").Value = 120481
$excel.CutCopyMode = $false

# Column widths (closest achievable increments; engine quantizes to 1/6 character units)
$intl.Columns.Item(8).ColumnWidth = 14.166666666666666
$intl.Columns.Item(9).ColumnWidth = 15.451822916666666

$intl.Range("H8").Select() | Out-Null
